# Apply the "PO Forecast" update to the workbook:
#  1. Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#     to "Weekly_PO_Qty".
#  2. Rename the "Requested quantity" header on the "Monthly Trend" sheet
#     to "Monthly_PO_Qty".
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing
#     the ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename header ---------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename header ------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet after "Monthly Trend" ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline/sheetPr settings used elsewhere in the workbook
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold / centered / bordered) used by the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data (rows 2-26)
$dates = @(45123.99999999999,45130.99999999999,45151.99999999999,45158.99999999999,45172.99999999999,45186.99999999999,45200.99999999999,45207.99999999999,45214.99999999999,45221.99999999999,45235.99999999999,45242.99999999999,45249.99999999999,45256.99999999999,45277.99999999999,45298.99999999999,45312.99999999999,45319.99999999999,45326.99999999999,45333.99999999999,45340.99999999999,45347.99999999999,45354.99999999999,45361.99999999999,45368.99999999999)
$forecast = @(35,52,104,121,155,189,223,241,258,275,309,326,343,360,412,463,498,515,532,549,566,583,600,617,635)
$lower = @(-211.7023266504851,-188.090583239742,-146.2514041096101,-144.424632179195,-94.09424451181813,-35.95596459945347,-14.36854128795266,-24.2330783433521,8.39858420620682,29.04215644670952,70.86967278962966,102.6086794272301,95.06168800707232,105.1769924946073,161.5560599960373,202.2242024584393,238.0010009164325,257.5877330533106,297.8846076262525,295.9054594327412,317.2224734509805,340.2754510910459,357.8326749970553,374.9781712797485,388.5916076296653)
$upper = @(270.9751280373401,304.5553081863507,349.1564938416538,363.7406560796225,399.9187396277546,431.2013522860864,477.3652020660953,464.9638265902126,508.4214054391298,520.0354578191016,531.7190427461443,562.9627597786349,591.8639738495654,599.811972307598,654.899484228986,703.4751738352262,759.5036250157854,743.5168936530084,764.3024886148546,801.6367200807214,811.5431789081223,839.7701283998335,821.8973503035106,865.8198876954895,870.1230282341028)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($r, 2).Value = $forecast[$i]
    $wsForecast.Cells.Item($r, 3).Value = $lower[$i]
    $wsForecast.Cells.Item($r, 4).Value = $upper[$i]
}

# Match the date-format styling used on the other sheets' date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A26").PasteSpecial(-4122)

Write-Output "PO Forecast sheet created with $($dates.Length) data rows"
